$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataset")

$row = 8

# Copy B7:O7 *values* (not a typed re-entry) down into B8:O8 so that the
# numeric-looking text ("0"/"1") and the date-looking text ("12/03/2018")
# stay literal strings, exactly like every other data row in this sheet,
# instead of being auto-converted into real numbers / a date serial the
# way a plain `.Value = "0"` assignment would.
$ws.Range("B7:O7").Copy()
$ws.Range("B8:O8").PasteSpecial(-4163)  # xlPasteValues

# New row's id is one past the last existing id (row 7 = 6).
$ws.Cells.Item($row, 1).Value = 7

# Match the bold / thin-border / center+top look used for A4:A7.
$idCell = $ws.Cells.Item($row, 1)
$idCell.Font.Bold = $true
$idCell.HorizontalAlignment = -4108  # xlCenter
$idCell.VerticalAlignment = -4160    # xlTop
$idCell.Borders.LineStyle = 1        # xlContinuous
$idCell.Borders.Weight = 2           # xlThin

$wb.Save()
